# Apply updated crypto price/volume figures (coinranking.com snapshot refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells that are plain decimals (e.g. "213.64") need a leading apostrophe so
# Excel stores them as text (matching the original "Price" column formatting)
# instead of auto-converting them to floating-point numbers.

$ws.Range('D2').Value = '27.814.95'
$ws.Range('E2').Value = '  +1.39%  '

$ws.Range('D3').Value = '1.650.22'
$ws.Range('E3').Value = '  -0.34%  '

$ws.Range('D5').Value = '''213.64'
$ws.Range('E5').Value = '  +0.17%  '

$ws.Range('D6').Value = '''0.534'
$ws.Range('E6').Value = '  -0.72%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = '''23.18'
$ws.Range('E8').Value = '  -1.23%  '

$ws.Range('E10').Value = '  +0.57%  '

$ws.Range('D11').Value = '''0.0891'

$ws.Range('D12').Value = '1.883.88'
$ws.Range('E12').Value = '  -0.31%  '

$ws.Range('D13').Value = '1.650.81'
$ws.Range('E13').Value = '  -0.32%  '

$ws.Range('E14').Value = '  -0.36%  '

$ws.Range('E15').Value = '  -0.28%  '

$ws.Range('E16').Value = '  -1.38%  '

$ws.Range('D17').Value = '27.777.70'
$ws.Range('E17').Value = '  +1.30%  '

$ws.Range('D18').Value = '''234.44'
$ws.Range('E18').Value = '  +2.32%  '

$ws.Range('E19').Value = '  +4.01%  '

$ws.Range('E21').Value = '  +0.04%  '

$ws.Range('E22').Value = '  -0.62%  '

$ws.Range('E23').Value = '  +8.04%  '

$ws.Range('E24').Value = '  -3.92%  '

$ws.Range('D25').Value = '''150.68'
$ws.Range('E25').Value = '  +2.55%  '

$ws.Range('D26').Value = '''6.98'
$ws.Range('E26').Value = '  -1.11%  '

$ws.Range('E27').Value = '  -1.55%  '

$ws.Range('D28').Value = '''15.71'
$ws.Range('E28').Value = '  +0.44%  '

$ws.Range('E30').Value = '  +0.57%  '

$ws.Range('E31').Value = '  -0.91%  '

$ws.Range('E32').Value = '  +0.66%  '

$ws.Range('E33').Value = '  +1.56%  '

$ws.Range('D34').Value = '1.445.60'
$ws.Range('E34').Value = '  +1.57%  '

$ws.Range('E35').Value = '  +1.82%  '

$ws.Range('E36').Value = '  -1.10%  '

$ws.Range('D37').Value = '''0.571'
$ws.Range('E37').Value = '  +0.81%  '

$ws.Range('E38').Value = '  -2.11%  '

$ws.Range('E39').Value = '  -0.16%  '

$ws.Range('D40').Value = '''0.879'
$ws.Range('E40').Value = '  +11.44%  '

$ws.Range('D41').Value = '''1.03'
$ws.Range('E41').Value = '  -1.00%  '

$ws.Range('E42').Value = '  +0.11%  '

$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('D44').Value = '''66.72'
$ws.Range('E44').Value = '  +2.30%  '

$ws.Range('D45').Value = '''2.46'
$ws.Range('E45').Value = '  -0.83%  '

$ws.Range('E46').Value = '  +2.13%  '

$ws.Range('D47').Value = '1.792.97'
$ws.Range('E47').Value = '  -0.32%  '

$ws.Range('E48').Value = '  +4.53%  '

$ws.Range('D49').Value = '''86.49'
$ws.Range('E49').Value = '  -1.55%  '

$ws.Range('D50').Value = '0.0₆0108'
$ws.Range('E50').Value = '  +3.03%  '

$ws.Range('E51').Value = '  -1.10%  '
